# Auto-generated update of Ultima_Profits market-price sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 137.46
$ws.Range("I15").Value = 137.46
$ws.Range("K15").Value = 412.38
$ws.Range("M15").Value = -243.38

$ws.Range("H32").Value = 1900.4
$ws.Range("J32").Value = 1900.4
$ws.Range("L32").Value = 1900.4
$ws.Range("N32").Value = -2552.4

$ws.Range("H62").Value = 1475
$ws.Range("I62").Value = 1650
$ws.Range("J62").Value = 1300
$ws.Range("K62").Value = 1650
$ws.Range("L62").Value = 1300
$ws.Range("M62").Value = -1026
$ws.Range("N62").Value = -2548

$ws.Range("H65").Value = 1475
$ws.Range("I65").Value = 1650
$ws.Range("J65").Value = 1300
$ws.Range("K65").Value = 8250
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = -5130
$ws.Range("N65").Value = -12740

$ws.Range("H133").Value = 31197.777
$ws.Range("J133").Value = 31197.777
$ws.Range("L133").Value = 31197.777
$ws.Range("N133").Value = -41317.777

$ws.Range("H137").Value = 5735.353
$ws.Range("I137").Value = 1191
$ws.Range("J137").Value = 14066.667
$ws.Range("K137").Value = 3573
$ws.Range("L137").Value = 42200.001
$ws.Range("M137").Value = -1023
$ws.Range("N137").Value = -47300.001

$ws.Range("H141").Value = 1169.6
$ws.Range("I141").Value = 945.2174
$ws.Range("J141").Value = 3750
$ws.Range("K141").Value = 2835.6522
$ws.Range("L141").Value = 11250
$ws.Range("M141").Value = 2344.3478
$ws.Range("N141").Value = -21610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13942.562
$ws.Range("I32").Value = 13487.594
$ws.Range("J32").Value = 19857.143
$ws.Range("K32").Value = 13487.594
$ws.Range("L32").Value = 19857.143
$ws.Range("M32").Value = -13200.594
$ws.Range("N32").Value = -20431.143

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1159.7693
$ws.Range("I22").Value = 1322.7
$ws.Range("J22").Value = 616.6667
$ws.Range("K22").Value = 1322.7
$ws.Range("L22").Value = 616.6667
$ws.Range("M22").Value = -1149.7
$ws.Range("N22").Value = -962.6667

$ws.Range("H36").Value = 1968.5
$ws.Range("I36").Value = 1968.5
$ws.Range("K36").Value = 1968.5
$ws.Range("M36").Value = -1434.5

$ws.Range("H134").Value = 3818.4075
$ws.Range("I134").Value = 2507.5
$ws.Range("J134").Value = 5725.1816
$ws.Range("K134").Value = 7522.5
$ws.Range("L134").Value = 17175.5448
$ws.Range("M134").Value = -4987.5
$ws.Range("N134").Value = -22245.5448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3483.1606
$ws.Range("I31").Value = 1662.6428
$ws.Range("J31").Value = 5303.6787
$ws.Range("K31").Value = 1662.6428
$ws.Range("L31").Value = 5303.6787
$ws.Range("M31").Value = -1367.6428
$ws.Range("N31").Value = -5893.6787

$ws.Range("H34").Value = 3483.1606
$ws.Range("I34").Value = 1662.6428
$ws.Range("J34").Value = 5303.6787
$ws.Range("K34").Value = 1662.6428
$ws.Range("L34").Value = 5303.6787
$ws.Range("M34").Value = -1460.6428
$ws.Range("N34").Value = -5707.6787

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4766.6665
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4766.6665
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 14299.9995
$ws.Range("N75").Value = -16295.9995

$ws.Range("H78").Value = 4766.6665
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4766.6665
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 42899.9985
$ws.Range("N78").Value = -52883.9985

$ws.Range("H122").Value = 1604.25
$ws.Range("I122").Value = 11004
$ws.Range("J122").Value = 1195.5652
$ws.Range("K122").Value = 99036
$ws.Range("L122").Value = 10760.0868
$ws.Range("M122").Value = -96586
$ws.Range("N122").Value = -15660.0868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4525
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 5050
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 15150
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -20050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9313.666999999999
$ws.Range("I7").Value = 6958.9
$ws.Range("J7").Value = 11454.363
$ws.Range("K7").Value = 6958.9
$ws.Range("L7").Value = 11454.363
$ws.Range("M7").Value = -6846.9
$ws.Range("N7").Value = -11678.363

$ws.Range("H22").Value = 1469.0714
$ws.Range("I22").Value = 870.2
$ws.Range("J22").Value = 1801.7778
$ws.Range("K22").Value = 870.2
$ws.Range("L22").Value = 1801.7778
$ws.Range("M22").Value = -575.2
$ws.Range("N22").Value = -2391.7778

$ws.Range("H27").Value = 1469.0714
$ws.Range("I27").Value = 870.2
$ws.Range("J27").Value = 1801.7778
$ws.Range("K27").Value = 870.2
$ws.Range("L27").Value = 1801.7778
$ws.Range("M27").Value = -763.2
$ws.Range("N27").Value = -2015.7778

$ws.Range("H40").Value = 3445.6428
$ws.Range("I40").Value = 8740.429
$ws.Range("J40").Value = 1680.7142
$ws.Range("K40").Value = 8740.429
$ws.Range("L40").Value = 1680.7142
$ws.Range("M40").Value = -8604.429
$ws.Range("N40").Value = -1952.7142

$ws.Range("H122").Value = 15077.556
$ws.Range("I122").Value = 50000
$ws.Range("J122").Value = 10712.25
$ws.Range("K122").Value = 150000
$ws.Range("L122").Value = 32136.75
$ws.Range("M122").Value = -147550
$ws.Range("N122").Value = -37036.75

$ws.Range("H126").Value = 9313.666999999999
$ws.Range("I126").Value = 6958.9
$ws.Range("J126").Value = 11454.363
$ws.Range("K126").Value = 20876.7
$ws.Range("L126").Value = 34363.089
$ws.Range("M126").Value = -18406.7
$ws.Range("N126").Value = -39303.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 32041.5
$ws.Range("J109").Value = 32041.5
$ws.Range("L109").Value = 32041.5
$ws.Range("N109").Value = -34815.5

$ws.Range("H122").Value = 2545
$ws.Range("I122").Value = 3317.8333
$ws.Range("J122").Value = 999.3333
$ws.Range("K122").Value = 9953.499899999999
$ws.Range("L122").Value = 2997.9999
$ws.Range("M122").Value = -7503.499899999999
$ws.Range("N122").Value = -7897.9999

$ws.Range("H126").Value = 2017.9269
$ws.Range("I126").Value = 2029.762
$ws.Range("K126").Value = 6089.286
$ws.Range("M126").Value = -3619.286

# ARM!N56 is removed entirely by the source refresh (row has no HQ-profit figure anymore)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N56").ClearContents()

# CUL!M75 and CUL!M78 are removed entirely (no NQ-profit figure for these rows anymore)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()
